$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3637.5386
$ws.Range("I32").Value = 2977.1428
$ws.Range("J32").Value = 4408
$ws.Range("K32").Value = 2977.1428
$ws.Range("L32").Value = 4408
$ws.Range("M32").Value = -2651.1428
$ws.Range("N32").Value = -5060
$ws.Range("H69").Value = 15465.767
$ws.Range("I69").Value = 12373.25
$ws.Range("J69").Value = 16590.318
$ws.Range("K69").Value = 37119.75
$ws.Range("L69").Value = 49770.954
$ws.Range("M69").Value = -36245.75
$ws.Range("N69").Value = -51518.954
$ws.Range("H70").Value = 3995.6
$ws.Range("I70").Value = 989
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 2967
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -2697
$ws.Range("N70").Value = -18540
$ws.Range("H72").Value = 15465.767
$ws.Range("I72").Value = 12373.25
$ws.Range("J72").Value = 16590.318
$ws.Range("K72").Value = 111359.25
$ws.Range("L72").Value = 149312.862
$ws.Range("M72").Value = -106991.25
$ws.Range("N72").Value = -158048.862
$ws.Range("H73").Value = 3995.6
$ws.Range("I73").Value = 989
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 2967
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -2031
$ws.Range("N73").Value = -19872
$ws.Range("H113").Value = 2238.4285
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 2311.5
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2311.5
$ws.Range("M113").Value = 1454
$ws.Range("N113").Value = -8819.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4572.985
$ws.Range("I32").Value = 2372.1406
$ws.Range("J32").Value = 75000
$ws.Range("K32").Value = 2372.1406
$ws.Range("L32").Value = 75000
$ws.Range("M32").Value = -2085.1406
$ws.Range("N32").Value = -75574
$ws.Range("H45").Value = 4297
$ws.Range("I45").Value = 3641.5
$ws.Range("J45").Value = 4559.2
$ws.Range("K45").Value = 3641.5
$ws.Range("L45").Value = 4559.2
$ws.Range("M45").Value = -3264.5
$ws.Range("N45").Value = -5313.2
$ws.Range("H61").Value = 4188.737
$ws.Range("I61").Value = 2260.4285
$ws.Range("J61").Value = 9588
$ws.Range("K61").Value = 2260.4285
$ws.Range("L61").Value = 9588
$ws.Range("M61").Value = -2048.4285
$ws.Range("N61").Value = -10012
$ws.Range("H74").Value = 3293.5386
$ws.Range("I74").Value = 2159.6667
$ws.Range("J74").Value = 16900
$ws.Range("K74").Value = 2159.6667
$ws.Range("L74").Value = 16900
$ws.Range("M74").Value = -1285.6667
$ws.Range("N74").Value = -18648
$ws.Range("H77").Value = 3293.5386
$ws.Range("I77").Value = 2159.6667
$ws.Range("J77").Value = 16900
$ws.Range("K77").Value = 10798.3335
$ws.Range("L77").Value = 84500
$ws.Range("M77").Value = -6430.333500000001
$ws.Range("N77").Value = -93236
$ws.Range("H122").Value = 2924.56
$ws.Range("I122").Value = 2130.2727
$ws.Range("J122").Value = 8749.333000000001
$ws.Range("K122").Value = 6390.8181
$ws.Range("L122").Value = 26247.999
$ws.Range("M122").Value = -3940.8181
$ws.Range("N122").Value = -31147.999
$ws.Range("H136").Value = 4188.737
$ws.Range("I136").Value = 2260.4285
$ws.Range("J136").Value = 9588
$ws.Range("K136").Value = 6781.2855
$ws.Range("L136").Value = 28764
$ws.Range("M136").Value = -4231.2855
$ws.Range("N136").Value = -33864

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H64").Value = 794.5
$ws.Range("I64").Value = 794.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 794.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -569.5
$ws.Range("H67").Value = 794.5
$ws.Range("I67").Value = 794.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 794.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -14.5
$ws.Range("H130").Value = 80000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 80000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H134").Value = 4237.186
$ws.Range("I134").Value = 2072.5625
$ws.Range("J134").Value = 10534.272
$ws.Range("K134").Value = 6217.6875
$ws.Range("L134").Value = 31602.816
$ws.Range("M134").Value = -3682.6875
$ws.Range("N134").Value = -36672.81600000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10534.069
$ws.Range("I31").Value = 4417.25
$ws.Range("J31").Value = 14851.823
$ws.Range("K31").Value = 4417.25
$ws.Range("L31").Value = 14851.823
$ws.Range("M31").Value = -4122.25
$ws.Range("N31").Value = -15441.823
$ws.Range("H34").Value = 10534.069
$ws.Range("I34").Value = 4417.25
$ws.Range("J34").Value = 14851.823
$ws.Range("K34").Value = 4417.25
$ws.Range("L34").Value = 14851.823
$ws.Range("M34").Value = -4215.25
$ws.Range("N34").Value = -15255.823
$ws.Range("H86").Value = 4120.6665
$ws.Range("I86").Value = 3860.75
$ws.Range("J86").Value = 6200
$ws.Range("K86").Value = 3860.75
$ws.Range("L86").Value = 6200
$ws.Range("M86").Value = -2737.75
$ws.Range("N86").Value = -8446
$ws.Range("H89").Value = 4120.6665
$ws.Range("I89").Value = 3860.75
$ws.Range("J89").Value = 6200
$ws.Range("K89").Value = 19303.75
$ws.Range("L89").Value = 31000
$ws.Range("M89").Value = -13687.75
$ws.Range("N89").Value = -42232
$ws.Range("H132").Value = 2516.1133
$ws.Range("I132").Value = 1795.4348
$ws.Range("J132").Value = 7252
$ws.Range("K132").Value = 5386.3044
$ws.Range("L132").Value = 21756
$ws.Range("M132").Value = -2856.3044
$ws.Range("N132").Value = -26816
$ws.Range("H134").Value = 9441.166999999999
$ws.Range("I134").Value = 5484.706
$ws.Range("J134").Value = 14615
$ws.Range("K134").Value = 16454.118
$ws.Range("L134").Value = 43845
$ws.Range("M134").Value = -13919.118
$ws.Range("N134").Value = -48915
$ws.Range("H141").Value = 325453.8
$ws.Range("I141").Value = 57432
$ws.Range("J141").Value = 425962
$ws.Range("K141").Value = 57432
$ws.Range("L141").Value = 425962
$ws.Range("M141").Value = -52252
$ws.Range("N141").Value = -436322

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 533.6667
$ws.Range("I8").Value = 533.6667
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1601.0001
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1462.0001
$ws.Range("H16").Value = 7777
$ws.Range("I16").Value = 7777
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 23331
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -23158
$ws.Range("H106").Value = 16665.666
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 16665.666
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 49996.99800000001
$ws.Range("N106").Value = -51888.99800000001
$ws.Range("H140").Value = 4777.675
$ws.Range("I140").Value = 27640
$ws.Range("J140").Value = 2237.4167
$ws.Range("K140").Value = 82920
$ws.Range("L140").Value = 6712.250100000001
$ws.Range("M140").Value = -77740
$ws.Range("N140").Value = -17072.2501

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8332.454
$ws.Range("I70").Value = 7522.5713
$ws.Range("J70").Value = 9749.75
$ws.Range("K70").Value = 7522.5713
$ws.Range("L70").Value = 9749.75
$ws.Range("M70").Value = -7252.5713
$ws.Range("N70").Value = -10289.75
$ws.Range("H73").Value = 8332.454
$ws.Range("I73").Value = 7522.5713
$ws.Range("J73").Value = 9749.75
$ws.Range("K73").Value = 7522.5713
$ws.Range("L73").Value = 9749.75
$ws.Range("M73").Value = -6586.5713
$ws.Range("N73").Value = -11621.75
$ws.Range("H97").Value = 708.0769
$ws.Range("I97").Value = 419.4
$ws.Range("J97").Value = 1670.3334
$ws.Range("K97").Value = 419.4
$ws.Range("L97").Value = 1670.3334
$ws.Range("M97").Value = 76.60000000000002
$ws.Range("N97").Value = -2662.3334
$ws.Range("H122").Value = 10787.889
$ws.Range("I122").Value = 6673.25
$ws.Range("J122").Value = 14079.6
$ws.Range("K122").Value = 20019.75
$ws.Range("L122").Value = 42238.8
$ws.Range("M122").Value = -17569.75
$ws.Range("N122").Value = -47138.8
$ws.Range("H126").Value = 6540.25
$ws.Range("I126").Value = 6540.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 19620.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17150.75
$ws.Range("H132").Value = 2878.4146
$ws.Range("I132").Value = 2613.5405
$ws.Range("J132").Value = 5328.5
$ws.Range("K132").Value = 7840.6215
$ws.Range("L132").Value = 15985.5
$ws.Range("M132").Value = -5310.6215
$ws.Range("N132").Value = -21045.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6157.4443
$ws.Range("I68").Value = 2780
$ws.Range("J68").Value = 10379.25
$ws.Range("K68").Value = 2780
$ws.Range("L68").Value = 10379.25
$ws.Range("M68").Value = -2031
$ws.Range("N68").Value = -11877.25
$ws.Range("H71").Value = 6157.4443
$ws.Range("I71").Value = 2780
$ws.Range("J71").Value = 10379.25
$ws.Range("K71").Value = 13900
$ws.Range("L71").Value = 51896.25
$ws.Range("M71").Value = -10156
$ws.Range("N71").Value = -59384.25
$ws.Range("H100").Value = 4166.2
$ws.Range("I100").Value = 3357.1538
$ws.Range("J100").Value = 5668.7144
$ws.Range("K100").Value = 3357.1538
$ws.Range("L100").Value = 5668.7144
$ws.Range("M100").Value = -2816.1538
$ws.Range("N100").Value = -6750.7144
$ws.Range("H122").Value = 5008.75
$ws.Range("I122").Value = 3930.0435
$ws.Range("J122").Value = 9970.799999999999
$ws.Range("K122").Value = 11790.1305
$ws.Range("L122").Value = 29912.4
$ws.Range("M122").Value = -9340.130500000001
$ws.Range("N122").Value = -34812.39999999999
$ws.Range("H132").Value = 8461.048000000001
$ws.Range("I132").Value = 8442.444
$ws.Range("J132").Value = 8475
$ws.Range("K132").Value = 25327.332
$ws.Range("L132").Value = 25425
$ws.Range("M132").Value = -22797.332
$ws.Range("N132").Value = -30485
$ws.Range("H136").Value = 8522.93
$ws.Range("I136").Value = 4453.579
$ws.Range("J136").Value = 10009.808
$ws.Range("K136").Value = 13360.737
$ws.Range("L136").Value = 30029.424
$ws.Range("M136").Value = -10810.737
$ws.Range("N136").Value = -35129.424

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 45655.832
$ws.Range("I81").Value = 86248.414
$ws.Range("J81").Value = 5063.25
$ws.Range("K81").Value = 172496.828
$ws.Range("L81").Value = 10126.5
$ws.Range("M81").Value = -171435.828
$ws.Range("N81").Value = -12248.5
$ws.Range("H84").Value = 45655.832
$ws.Range("I84").Value = 86248.414
$ws.Range("J84").Value = 5063.25
$ws.Range("K84").Value = 862484.14
$ws.Range("L84").Value = 50632.5
$ws.Range("M84").Value = -857180.14
$ws.Range("N84").Value = -61240.5
$ws.Range("H107").Value = 1816.138
$ws.Range("I107").Value = 2030.1765
$ws.Range("J107").Value = 1512.9166
$ws.Range("K107").Value = 6090.529500000001
$ws.Range("L107").Value = 4538.7498
$ws.Range("M107").Value = -4170.529500000001
$ws.Range("N107").Value = -8378.7498
$ws.Range("H132").Value = 2782.5312
$ws.Range("I132").Value = 1608.6428
$ws.Range("J132").Value = 10999.75
$ws.Range("K132").Value = 4825.928400000001
$ws.Range("L132").Value = 32999.25
$ws.Range("M132").Value = -2295.928400000001
$ws.Range("N132").Value = -38059.25
